$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add basic data references" - new column O holding a reference/flag
# column (Coarse_seds_subsurface) with a 0 default for every data row.
$ws.Range("O1").Value = "Coarse_seds_subsurface"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
}

# "... and flux box outline" - the header row grows taller to show the
# wrapped/outlined header box across the now-wider table.
$ws.Rows.Item(1).RowHeight = 58

# Selection left on the sheet when the workbook was last saved.
$ws.Range("H20").Select()
